# Update workbook metadata and element-table content to reflect the
# LinuxForHealth rebrand (formerly Alvearie / ibm.com) and related
# version/date bump.

$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet -------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/age-group"
$wsMeta.Range("B3").Value = "8.0.0"
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# --- "Elements" sheet ---------------------------------------------------
$wsElements = $wb.Worksheets.Item("Elements")

# The base "Extension" row no longer carries the ele-1/ext-1 constraint
# text in its Constraint(s) column (it now only appears on the
# Extension.extension row below).
$wsElements.Range("AI2").Value = ""

# Fixed Value / Binding Value Set URLs also move from ibm.com to
# linuxforhealth.org.
$wsElements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/age-group"
$wsElements.Range("Y7").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/age-group"

# The "Binding Value Set" column widened to fit the longer URL text.
$wsElements.Columns.Item(25).ColumnWidth = 49.3
